# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mutates the "before" localization-status workbook so that
# it reflects a completed handback cycle:
#   * Status flips from "Ready for handoff" to "Handed back: in sync
#     with en-US" everywhere it appears (Overview + per-locale sheets).
#   * The per-locale tables (zh-cn, de-de) gain a populated "Latest
#     Target File" / "Latest Handback File" / "Latest Handback
#     DateTime" for rows 2 & 3, with the Target File cell turned into a
#     hyperlink back to the source markdown file.
#   * A handful of columns are widened so the new, longer values aren't
#     clipped.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) == FF6495ED, matches the workbook's existing HyperLink style

# -------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell that currently carries the old status)
# -------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# -------------------------------------------------------------------
# 2. Widen columns that now hold longer strings.
# -------------------------------------------------------------------
# Overview: "zh-cn" / "de-de" status columns (E, F)
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# Per-locale sheets: Status column (C) and the new Target/Handback File columns (I, J)
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.144371396019366
    $ws.Columns.Item(9).ColumnWidth  = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

# -------------------------------------------------------------------
# 3. zh-cn sheet: populate "Latest Target File" (I) / "Latest Handback
#    File" (J) / "Latest Handback DateTime" (K) for rows 2 & 3.
# -------------------------------------------------------------------
$mdTargetName = "93808e0a-e246-4825-aff9-e47cfeec904e.md"
$mdTargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/170a347512d49373f517d40da7e312d663b77d47/e2e/93808e0a-e246-4825-aff9-e47cfeec904e.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdTargetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdTargetName)
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("J2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 15:04:26"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdTargetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdTargetName)
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor
$wsZhCn.Range("J3").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-20 15:04:26"

# -------------------------------------------------------------------
# 4. de-de sheet: same shape, different target file + later timestamp.
# -------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdTargetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdTargetName)
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("J2").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 15:04:32"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdTargetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdTargetName)
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor
$wsDeDe.Range("J3").Value = "93808e0a-e246-4825-aff9-e47cfeec904e.44997de73ba84e5487eb6facad0287e9ef16634f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-20 15:04:32"
